# Adjust the row heights of the header/data rows (1-18) on the active sheet
# from 19.5 to 18.75, matching the height already used by the rest of the
# sheet's rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:18").RowHeight = 18.75
